# "I've learned more topics"
#
# This sheet is a simple two-column DSA topic tracker (A: Topic, B: Status).
# The edit:
#   - Renames the old row 12 topic ("HashMap or HashTable") to "Queues"
#     and marks it Done (it previously had no status).
#   - Adds three more topic rows: HashMap (Done), Tree (no status yet),
#     plus two trailing blank rows that are already formatted for future entries.
#   - Leaves the scroll position/selection near the newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the existing row formatting (borders/fills) down into the four new
# rows before we touch any values, so the new cells look like the rest of
# the table (plain bordered cell in col A, light-green fill in col B).
$ws.Range("A12:B12").Copy()
$ws.Range("A13:B16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 12 used to just be "HashMap or HashTable" with an empty status; it is
# now "Queues" and marked Done.
$ws.Range("A12").Value = "Queues"
$ws.Range("B12").Value = "Done"

# New row 13: HashMap, Done.
$ws.Range("A13").Value = "HashMap"
$ws.Range("B13").Value = "Done"

# New row 14: Tree, status not filled in yet.
$ws.Range("A14").Value = "Tree"

# Rows 15-16 stay blank (just formatted), ready for future topics.

# Match the author's final view state: scrolled down a bit with A15 selected.
$ws.Range("A4").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("A15").Select() | Out-Null
